$d = $word.ActiveDocument

# The document ends with a single empty paragraph that only holds the
# "_GoBack" bookmark. The edit:
#   1. inserts two new strikethrough to-do paragraphs before it,
#   2. strips the <w:rFonts cstheme="minorHAnsi"/> from that bookmark
#      paragraph's mark run-properties (keeps <w:strike/>),
#   3. appends a new empty paragraph (rFonts + strike) after it.
#
# Replace the whole final paragraph's Range in one shot with the four
# target paragraphs so every run/paragraph-mark property ends up exactly
# as specified (InsertXML replaces the contents of the Range it is
# called on).

$count = $d.Paragraphs.Count
$last = $d.Paragraphs.Item($count)
$r = $last.Range

$xml = @"
<pkg:xmlData xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:p><w:pPr><w:rPr><w:strike/></w:rPr></w:pPr><w:r><w:rPr><w:strike/></w:rPr><w:t>Estudar ANOVA – interpretação de resultados (</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:strike/></w:rPr><w:t>lm</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:strike/></w:rPr><w:t xml:space="preserve"> e </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:strike/></w:rPr><w:t>aov</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:strike/></w:rPr><w:t>)</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:strike/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:strike/></w:rPr><w:t xml:space="preserve">Normalidade e </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:strike/></w:rPr><w:t>Levene</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:strike/></w:rPr><w:t xml:space="preserve"> para cada combinação</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:strike/></w:rPr></w:pPr><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p><w:p><w:pPr><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:strike/></w:rPr></w:pPr></w:p></pkg:xmlData>
"@

$r.InsertXML($xml)
Write-Output "paragraphs now: $($d.Paragraphs.Count)"
